# Scheduled market-data refresh: recompute currentAveragePrice* / LevePrice* /
# LeveProfit* (columns H-N) for the leves whose Universalis snapshot changed,
# across all 8 crafter-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 787.13043
$ws.Range("I19").Value = 783.8461
$ws.Range("J19").Value = 791.4
$ws.Range("K19").Value = 783.8461
$ws.Range("L19").Value = 791.4
$ws.Range("M19").Value = -608.8461
$ws.Range("N19").Value = -1141.4

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 2296.9033
$ws.Range("I116").Value = 1664.091
$ws.Range("J116").Value = 2644.95
$ws.Range("K116").Value = 1664.091
$ws.Range("L116").Value = 2644.95
$ws.Range("M116").Value = 1777.909
$ws.Range("N116").Value = -9528.950000000001

# Row 135: For Tired Minds / Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 26318070
$ws.Range("J135").Value = 71433060
$ws.Range("L135").Value = 642897540
$ws.Range("N135").Value = -642902610

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2271.8
$ws.Range("I138").Value = 1594.85
$ws.Range("J138").Value = 2813.36
$ws.Range("K138").Value = 4784.549999999999
$ws.Range("L138").Value = 8440.08
$ws.Range("M138").Value = 355.4500000000007
$ws.Range("N138").Value = -18720.08


$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 3659
$ws.Range("I74").Value = 1001.3333
$ws.Range("J74").Value = 4944.968
$ws.Range("K74").Value = 1001.3333
$ws.Range("L74").Value = 4944.968
$ws.Range("M74").Value = -127.3333
$ws.Range("N74").Value = -6692.968

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 3659
$ws.Range("I77").Value = 1001.3333
$ws.Range("J77").Value = 4944.968
$ws.Range("K77").Value = 5006.6665
$ws.Range("L77").Value = 24724.84
$ws.Range("M77").Value = -638.6665000000003
$ws.Range("N77").Value = -33460.84

# Row 86: Sir, Dost Thou Even Heft / Adamantite Chain Hose of Fending
$ws.Range("H86").Value = 23650
$ws.Range("I86").Value = 1300
$ws.Range("J86").Value = 46000
$ws.Range("K86").Value = 1300
$ws.Range("L86").Value = 46000
$ws.Range("M86").Value = -114
$ws.Range("N86").Value = -48372

# Row 89: Men in Adamantite (L) / Adamantite Chain Hose of Fending
$ws.Range("H89").Value = 23650
$ws.Range("I89").Value = 1300
$ws.Range("J89").Value = 46000
$ws.Range("K89").Value = 3900
$ws.Range("L89").Value = 138000
$ws.Range("M89").Value = 2028
$ws.Range("N89").Value = -149856

# Row 119: Trial and Error / Dwarven Mythril Chainmail of Fending
$ws.Range("H119").Value = 25256.75
$ws.Range("J119").Value = 25256.75
$ws.Range("L119").Value = 25256.75
$ws.Range("N119").Value = -34932.75


$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight / Mythrite Nugget
$ws.Range("H64").Value = 1511.75
$ws.Range("I64").Value = 3325.4285
$ws.Range("J64").Value = 764.94116
$ws.Range("K64").Value = 3325.4285
$ws.Range("L64").Value = 764.94116
$ws.Range("M64").Value = -3100.4285
$ws.Range("N64").Value = -1214.94116

# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Range("H67").Value = 1511.75
$ws.Range("I67").Value = 3325.4285
$ws.Range("J67").Value = 764.94116
$ws.Range("K67").Value = 3325.4285
$ws.Range("L67").Value = 764.94116
$ws.Range("M67").Value = -2545.4285
$ws.Range("N67").Value = -2324.94116

# Row 80: Unbreaker / Titanium Ingot
$ws.Range("H80").Value = 750.7857
$ws.Range("I80").Value = 121.5
$ws.Range("J80").Value = 1222.75
$ws.Range("K80").Value = 121.5
$ws.Range("L80").Value = 1222.75
$ws.Range("M80").Value = 876.5
$ws.Range("N80").Value = -3218.75

# Row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Range("H83").Value = 750.7857
$ws.Range("I83").Value = 121.5
$ws.Range("J83").Value = 1222.75
$ws.Range("K83").Value = 607.5
$ws.Range("L83").Value = 6113.75
$ws.Range("M83").Value = 4384.5
$ws.Range("N83").Value = -16097.75


$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 1729.3043
$ws.Range("I31").Value = 984.93335
$ws.Range("J31").Value = 3125
$ws.Range("K31").Value = 984.93335
$ws.Range("L31").Value = 3125
$ws.Range("M31").Value = -689.93335
$ws.Range("N31").Value = -3715

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 1729.3043
$ws.Range("I34").Value = 984.93335
$ws.Range("J34").Value = 3125
$ws.Range("K34").Value = 984.93335
$ws.Range("L34").Value = 3125
$ws.Range("M34").Value = -782.93335
$ws.Range("N34").Value = -3529

# Row 107: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 2778231
$ws.Range("I107").Value = 3788301.8
$ws.Range("K107").Value = 3788301.8
$ws.Range("M107").Value = -3786381.8

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 7214.857
$ws.Range("I132").Value = 13756
$ws.Range("J132").Value = 4598.4
$ws.Range("K132").Value = 41268
$ws.Range("L132").Value = 13795.2
$ws.Range("M132").Value = -38738
$ws.Range("N132").Value = -18855.2

# Row 138: Bow Out / Acacia Longbow
$ws.Range("H138").Value = 51290
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 51290
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 51290
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -61570


$ws = $wb.Worksheets.Item("CUL")
# Row 44: No More Dumpster Diving / Knight's Bread
$ws.Range("H44").Value = 1272.5
$ws.Range("I44").Value = 487.85715
$ws.Range("J44").Value = 2057.1428
$ws.Range("K44").Value = 1463.57145
$ws.Range("L44").Value = 6171.428400000001
$ws.Range("M44").Value = -1065.57145
$ws.Range("N44").Value = -6967.428400000001

# Row 97: The Frier Never Lies / Cottonseed Oil
$ws.Range("H97").Value = 2214.8572
$ws.Range("I97").Value = 400
$ws.Range("J97").Value = 2940.8
$ws.Range("K97").Value = 1200
$ws.Range("L97").Value = 8822.400000000001
$ws.Range("M97").Value = -704
$ws.Range("N97").Value = -9814.400000000001

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 11905926
$ws.Range("I122").Value = 18519008
$ws.Range("K122").Value = 166671072
$ws.Range("M122").Value = -166668622

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1335621.5
$ws.Range("I131").Value = 8857.143
$ws.Range("J131").Value = 1472200.2
$ws.Range("K131").Value = 26571.429
$ws.Range("L131").Value = 4416600.6
$ws.Range("M131").Value = -21531.429
$ws.Range("N131").Value = -4426680.6

# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 47619920
$ws.Range("I132").Value = 66667280
$ws.Range("K132").Value = 600005520
$ws.Range("M132").Value = -600002990


$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 1377.8235
$ws.Range("I102").Value = 1088.625
$ws.Range("J102").Value = 1634.8889
$ws.Range("K102").Value = 1088.625
$ws.Range("L102").Value = 1634.8889
$ws.Range("M102").Value = 533.375
$ws.Range("N102").Value = -4878.8889


$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 610.9394
$ws.Range("I61").Value = 520.5925999999999
$ws.Range("J61").Value = 1017.5
$ws.Range("K61").Value = 520.5925999999999
$ws.Range("L61").Value = 1017.5
$ws.Range("M61").Value = -318.5925999999999
$ws.Range("N61").Value = -1421.5

# Row 81: I Need Your Glove Tonight / Dragonskin Gloves of Healing
$ws.Range("H81").Value = 35000
$ws.Range("J81").Value = 35000
$ws.Range("L81").Value = 35000
$ws.Range("N81").Value = -36996

# Row 84: Halonic Drake Handlers (L) / Dragonskin Gloves of Healing
$ws.Range("H84").Value = 35000
$ws.Range("J84").Value = 35000
$ws.Range("L84").Value = 105000
$ws.Range("N84").Value = -114984

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 610.9394
$ws.Range("I113").Value = 520.5925999999999
$ws.Range("J113").Value = 1017.5
$ws.Range("K113").Value = 520.5925999999999
$ws.Range("L113").Value = 1017.5
$ws.Range("M113").Value = 1649.4074
$ws.Range("N113").Value = -5357.5

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 76620.664
$ws.Range("I132").Value = 125734.664
$ws.Range("J132").Value = 2949.6667
$ws.Range("K132").Value = 377203.992
$ws.Range("L132").Value = 8849.000100000001
$ws.Range("M132").Value = -374673.992
$ws.Range("N132").Value = -13909.0001


$ws = $wb.Worksheets.Item("WVR")
# Row 119: A Job Well Done / Dwarven Cotton Gaskins of Fending
$ws.Range("H119").Value = 26500
$ws.Range("J119").Value = 26500
$ws.Range("L119").Value = 26500
$ws.Range("N119").Value = -36176
